$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-25 Wednesday" "2024-09-26 Thursday"

Replace-Text "54×51=" "46×31="
Replace-Text "54×37=" "55×38="
Replace-Text "33×85=" "22×28="
Replace-Text "70×65=" "34×68="
Replace-Text "35×20=" "27×52="
Replace-Text "23×18=" "15×30="
Replace-Text "52×13=" "73×99="
Replace-Text "52×99=" "81×88="
Replace-Text "71×28=" "40×66="
Replace-Text "17×55=" "37×18="
Replace-Text "22×74=" "41×97="
Replace-Text "65×70=" "52×23="
Replace-Text "91×12=" "51×68="
Replace-Text "66×76=" "13×83="
Replace-Text "16×18=" "36×51="
Replace-Text "65×36=" "94×93="
Replace-Text "19×90=" "78×77="
Replace-Text "48×31=" "92×29="
Replace-Text "16×94=" "85×69="
Replace-Text "60×12=" "72×98="
Replace-Text "47×21=" "66×81="
Replace-Text "78×89=" "41×56="
Replace-Text "44×17=" "64×16="
Replace-Text "79×59=" "97×53="
Replace-Text "85×91=" "24×59="
